$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($addr, $val) {
    $r = $ws.Range($addr)
    $origStyle = $r.Style
    $r.NumberFormat = "@"
    $r.Value = $val
    $r.Style = $origStyle
}

Set-TextValue "D2" "43.124.56"
Set-TextValue "E2" "  +1.34%  "

Set-TextValue "D3" "2.377.27"
Set-TextValue "E3" "  +3.64%  "

Set-TextValue "E4" "  +0.11%  "

Set-TextValue "D5" "302.99"
Set-TextValue "E5" "  +0.80%  "

Set-TextValue "D6" "96.96"
Set-TextValue "E6" "  +1.84%  "

Set-TextValue "D7" "0.505"
Set-TextValue "E7" "  -0.08%  "

Set-TextValue "E8" "  -0.04%  "

Set-TextValue "D9" "0.499"
Set-TextValue "E9" "  +1.56%  "

Set-TextValue "D10" "34.23"
Set-TextValue "E10" "  -0.60%  "

Set-TextValue "D11" "0.0786"
Set-TextValue "E11" "  +0.87%  "

Set-TextValue "E12" "  +1.83%  "

Set-TextValue "D13" "18.32"
Set-TextValue "E13" "  -3.63%  "

Set-TextValue "E14" "  +1.85%  "

Set-TextValue "D15" "2.750.40"
Set-TextValue "E15" "  +3.81%  "

Set-TextValue "D16" "2.351.95"
Set-TextValue "E16" "  +1.04%  "

Set-TextValue "D17" "0.806"
Set-TextValue "E17" "  +3.86%  "

Set-TextValue "D18" "43.150.38"
Set-TextValue "E18" "  +1.60%  "

Set-TextValue "D19" "12.15"
Set-TextValue "E19" "  -0.52%  "

Set-TextValue "E20" "  +4.88%  "

Set-TextValue "D21" "0.0₃0887"
Set-TextValue "E21" "  +0.11%  "

Set-TextValue "D22" "68.38"
Set-TextValue "E22" "  +1.29%  "

Set-TextValue "D23" "235.26"

Set-TextValue "D24" "2.20"
Set-TextValue "E24" "  -3.09%  "

Set-TextValue "D25" "2.44"
Set-TextValue "E25" "  +2.12%  "

Set-TextValue "E26" "  -0.07%  "

Set-TextValue "D27" "24.85"
Set-TextValue "E27" "  +2.65%  "

Set-TextValue "E28" "  +0.38%  "

Set-TextValue "D29" "9.14"
Set-TextValue "E29" "  +1.34%  "

Set-TextValue "D30" "31.35"
Set-TextValue "E30" "  -0.89%  "

Set-TextValue "D32" "5.07"
Set-TextValue "E32" "  +1.89%  "

Set-TextValue "D33" "0.0745"
Set-TextValue "E33" "  +7.87%  "

Set-TextValue "D34" "17.37"
Set-TextValue "E34" "  -0.72%  "

Set-TextValue "E35" "  +5.43%  "

Set-TextValue "E36" "  +6.32%  "

Set-TextValue "B37" "WEMIXToken"
Set-TextValue "C37" "https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix"
Set-TextValue "D37" "2.32"
Set-TextValue "E37" "  -0.69%  "

Set-TextValue "B38" "RenderToken"
Set-TextValue "C38" "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
Set-TextValue "D38" "4.29"
Set-TextValue "E38" "  -2.11%  "

Set-TextValue "E39" "  +4.36%  "

Set-TextValue "D40" "22.27"
Set-TextValue "E40" "  +11.48%  "

Set-TextValue "E41" "  +0.07%  "

Set-TextValue "B42" "Maker"
Set-TextValue "C42" "https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr"
Set-TextValue "D42" "1.961.76"
Set-TextValue "E42" "  +0.84%  "

Set-TextValue "B43" "Monero"
Set-TextValue "C43" "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
Set-TextValue "D43" "103.21"
Set-TextValue "E43" "  -37.32%  "

Set-TextValue "E44" "  +1.11%  "

Set-TextValue "E45" "  +2.31%  "

Set-TextValue "D46" "2.74"
Set-TextValue "E46" "  +0.23%  "

Set-TextValue "D47" "9.12"
Set-TextValue "E47" "  -11.21%  "

Set-TextValue "D48" "2.612.04"
Set-TextValue "E48" "  +3.70%  "

Set-TextValue "D49" "52.61"
Set-TextValue "E49" "  -0.49%  "

Set-TextValue "E50" "  +2.26%  "

Set-TextValue "D51" "71.94"
Set-TextValue "E51" "  +2.15%  "
